$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop trailing rows 25-26 (their content no longer exists) ---
$ws.Rows.Item(26).Delete() | Out-Null
$ws.Rows.Item(25).Delete() | Out-Null

# --- Remove cells that existed before but are absent from the new layout ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()
$ws.Range("A24").Clear()

# --- Row 13 B/C needs literal text "01/01/2012"; re-use the formatting/value
#     already on B8/C8 (same text) via Copy so Excel does not reinterpret it as a date ---
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("C8").Copy($ws.Range("C13"))

# --- Set remaining cell values (rows 10-24; rows 1-9 stay untouched) ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C15").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("C18").Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios."
$ws.Range("C19").Value = "Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "WILLIAMS, D. B.; CARTER, C. B., Transmission Electron Microscopy: A Textbook for Materials Science, Springer, 2009.`nBOZZOLA, J. J.; RUSSELL, L. D. Electron Microscopy, Boston, Jones & Bartlett, 1999.`nHUNTER, E. Practical Electron Microscopy, Cambridge University Press, 1993.`nREIMER, L.; KOHL, H., Transmission Electron Microscopy: Physics of Image Formation, Springer, 2008.`nEATON, P.; WEST, P. Atomic Force Microscopy, Oxford University Press, 2010.`nMORITA, S.; WIESENDANGER, R.; MEYER, E. Noncontact Atomic Force Microscopy, Springer, 2002."
$ws.Range("C21").Value = "WILLIAMS, D. B.; CARTER, C. B., Transmission Electron Microscopy: A Textbook for Materials Science, Springer, 2009.`nBOZZOLA, J. J.; RUSSELL, L. D. Electron Microscopy, Boston, Jones & Bartlett, 1999.`nHUNTER, E. Practical Electron Microscopy, Cambridge University Press, 1993.`nREIMER, L.; KOHL, H., Transmission Electron Microscopy: Physics of Image Formation, Springer, 2008.`nEATON, P.; WEST, P. Atomic Force Microscopy, Oxford University Press, 2010.`nMORITA, S.; WIESENDANGER, R.; MEYER, E. Noncontact Atomic Force Microscopy, Springer, 2002."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOM3229 -  Métodos Experimentais da Física II  (Requisito)`n"
$ws.Range("C23").Value = "LOM3229 -  Métodos Experimentais da Física II  (Requisito)`n"
$ws.Range("B24").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Requisito)`n"
$ws.Range("C24").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Requisito)`n"

# --- Adjust row heights that differ from the sheet default ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(17).EntireRow.AutoFit()
$ws.Rows.Item(22).EntireRow.AutoFit()
